# rkap_detail.xlsx: make kd_aktifitas_rkap's "nama_aktifitas" (col C) match
# the more specific "uraian" (col D) value on every data row, instead of the
# old coarse category label. This also orphans several now-unused category
# strings (Suku Cadang, Bahan Pembantu, Lainnya - Bahan Bakar Kendaraan,
# Pemeliharaan - Kendaraan, Retribusi - Kendaraan, Retribusi - Parkir & Tol,
# Sewa - Kendaraan).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 23; $r++) {
    $ws.Range("C$r").Value = $ws.Range("D$r").Value2
}

# Column C now holds longer text ("nama_aktifitas" == "uraian"), widen it to fit.
$ws.Columns.Item(3).ColumnWidth = 29.7

# Update the view: scroll the window right a bit and select C8.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C8").Select()
